$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.040358665670715
$ws.Range("D2").Value = 1.04111634972258
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.039110467293263
$ws.Range("I2").Value = 1.0387933471537
$ws.Range("J2").Value = 1.045445701466658
$ws.Range("K2").Value = 1.043896775836738
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.041896586749309
$ws.Range("N2").Value = 1.046930354277038

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041590619980354
$ws.Range("D3").Value = 1.042030268635444
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.040970709757335
$ws.Range("I3").Value = 1.039140116245409
$ws.Range("J3").Value = 1.046321830691545
$ws.Range("K3").Value = 1.044620985184027
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.043564212872685
$ws.Range("N3").Value = 1.047807727705919

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.042386276866131
$ws.Range("D4").Value = 1.042620391068494
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.042172712801179
$ws.Range("I4").Value = 1.039362619725574
$ws.Range("J4").Value = 1.046886753365978
$ws.Range("K4").Value = 1.04508772796542
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.044641166917071
$ws.Range("N4").Value = 1.048373452635346

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.042720417226907
$ws.Range("D5").Value = 1.042868183818819
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.042677640996318
$ws.Range("I5").Value = 1.039455712058547
$ws.Range("J5").Value = 1.047123774147432
$ws.Range("K5").Value = 1.045283502309814
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.045093424263281
$ws.Range("N5").Value = 1.048610810013497

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.042776500270203
$ws.Range("D6").Value = 1.04290977212027
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.042762397931232
$ws.Range("I6").Value = 1.039471316443759
$ws.Range("J6").Value = 1.047163543377025
$ws.Range("K6").Value = 1.045316347719323
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.045169331614756
$ws.Range("N6").Value = 1.048650635719958

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.042390743049622
$ws.Range("D7").Value = 1.042623703241536
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.0421794612085
$ws.Range("I7").Value = 1.039363865388542
$ws.Range("J7").Value = 1.046889922301693
$ws.Range("K7").Value = 1.045090345653807
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.044647211929451
$ws.Range("N7").Value = 1.048376626071314

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.040775324165194
$ws.Range("D8").Value = 1.041425472424696
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.039739504222716
$ws.Range("I8").Value = 1.038910929591794
$ws.Range("J8").Value = 1.045742208252689
$ws.Range("K8").Value = 1.04414191468493
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.042460612890377
$ws.Range("N8").Value = 1.047227282136716

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.037917023504712
$ws.Range("D9").Value = 1.039304365040091
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.035426412199712
$ws.Range("I9").Value = 1.03809832510735
$ws.Range("J9").Value = 1.043704353483386
$ws.Range("K9").Value = 1.042456201099028
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.038590862230567
$ws.Range("N9").Value = 1.045186533379895

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.036003269222071
$ws.Range("D10").Value = 1.037883588598688
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.03254109103248
$ws.Range("I10").Value = 1.037546743356596
$ws.Range("J10").Value = 1.042335146332432
$ws.Range("K10").Value = 1.041322470804481
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.035999093572739
$ws.Range("N10").Value = 1.043815381797736

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035172566926915
$ws.Range("D11").Value = 1.037266740469992
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.031289182382354
$ws.Range("I11").Value = 1.037305541283467
$ws.Range("J11").Value = 1.041739683037891
$ws.Range("K11").Value = 1.040829154351283
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.0348738380289
$ws.Range("N11").Value = 1.043219072877017

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.034863694880323
$ws.Range("D12").Value = 1.03703736477211
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.030823769905493
$ws.Range("I12").Value = 1.03721559070134
$ws.Range("J12").Value = 1.04151810756901
$ws.Range("K12").Value = 1.040645549500897
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.034455403026015
$ws.Range("N12").Value = 1.042997182745554

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034929963218886
$ws.Range("D13").Value = 1.037086578054592
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.030923620714904
$ws.Range("I13").Value = 1.037234901616471
$ws.Range("J13").Value = 1.041565654166523
$ws.Range("K13").Value = 1.040684949953545
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.034545180020637
$ws.Range("N13").Value = 1.04304479686469

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035147041875115
$ws.Range("D14").Value = 1.03724778534393
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.031250719447835
$ws.Range("I14").Value = 1.037298113244369
$ws.Range("J14").Value = 1.041721375605609
$ws.Range("K14").Value = 1.040813984989515
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.034839259627486
$ws.Range("N14").Value = 1.043200739446081

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035280749627127
$ws.Range("D15").Value = 1.037347077114647
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.031452202588313
$ws.Range("I15").Value = 1.037337012587413
$ws.Range("J15").Value = 1.041817268375134
$ws.Range("K15").Value = 1.040893439226399
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.035020389745694
$ws.Range("N15").Value = 1.043296768394337

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.036058356777496
$ws.Range("D16").Value = 1.037924491825016
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.032624121141496
$ws.Range("I16").Value = 1.037562701151215
$ws.Range("J16").Value = 1.042374610315967
$ws.Range("K16").Value = 1.04135515960646
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.036073708633488
$ws.Range("N16").Value = 1.043854901824655

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.036545579959269
$ws.Range("D17").Value = 1.038286246205445
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.033358542328158
$ws.Range("I17").Value = 1.03770363522397
$ws.Range("J17").Value = 1.042723519816325
$ws.Range("K17").Value = 1.041644138225579
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.036733614096612
$ws.Range("N17").Value = 1.044204306816524

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.036829573192466
$ws.Range("D18").Value = 1.038497093210121
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.033786673221699
$ws.Range("I18").Value = 1.037785611794157
$ws.Range("J18").Value = 1.042926783512722
$ws.Range("K18").Value = 1.041812462911748
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.037118237024292
$ws.Range("N18").Value = 1.044407859170684

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.03692637453333
$ws.Range("D19").Value = 1.038568959877517
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.033932613775467
$ws.Range("I19").Value = 1.037813525101041
$ws.Range("J19").Value = 1.042996049047088
$ws.Range("K19").Value = 1.041869818108093
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.037249334941687
$ws.Range("N19").Value = 1.044477223070054

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.036493325825517
$ws.Range("D20").Value = 1.038247449778411
$ws.Range("E20").Value = 0.9894336180360677
$ws.Range("F20").Value = 1.033279771295495
$ws.Range("I20").Value = 1.037688537920705
$ws.Range("J20").Value = 1.042686110954875
$ws.Range("K20").Value = 1.041613157545669
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.036662842456854
$ws.Range("N20").Value = 1.0441668448302

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035083126293966
$ws.Range("D21").Value = 1.037200320772411
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.03115440813123
$ws.Range("I21").Value = 1.037279508879822
$ws.Range("J21").Value = 1.041675530437801
$ws.Range("K21").Value = 1.040775997500071
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.03475267346118
$ws.Range("N21").Value = 1.043154829172876

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.034194668098884
$ws.Range("D22").Value = 1.036540495890947
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.029815798303067
$ws.Range("I22").Value = 1.037020267168642
$ws.Range("J22").Value = 1.041037856860123
$ws.Range("K22").Value = 1.040247527059527
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.033548977133296
$ws.Range("N22").Value = 1.042516250025576

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.034665830374104
$ws.Range("D23").Value = 1.03689042067024
$ws.Range("E23").Value = 0.9879432794636464
$ws.Range("F23").Value = 1.030525644813585
$ws.Range("I23").Value = 1.037157892987555
$ws.Range("J23").Value = 1.041376117676024
$ws.Range("K23").Value = 1.040527881039524
$ws.Range("L23").Value = 0.9917760702887611
$ws.Range("M23").Value = 1.034187339668256
$ws.Range("N23").Value = 1.042854991210634

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.036516937823805
$ws.Range("D24").Value = 1.038264980702371
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.033315365288627
$ws.Range("I24").Value = 1.037695360446163
$ws.Range("J24").Value = 1.042703015179045
$ws.Range("K24").Value = 1.04162715709684
$ws.Range("L24").Value = 0.9929938892766441
$ws.Range("M24").Value = 1.036694822005404
$ws.Range("N24").Value = 1.044183773060307

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.038657386807605
$ws.Range("D25").Value = 1.039853888949609
$ws.Range("E25").Value = 0.9912096547607051
$ws.Range("F25").Value = 1.036543140532937
$ws.Range("I25").Value = 1.038310129730297
$ws.Range("J25").Value = 1.044233043540729
$ws.Range("K25").Value = 1.042893732277294
$ws.Range("L25").Value = 0.9944092447426416
$ws.Range("M25").Value = 1.039593336167115
$ws.Range("N25").Value = 1.045715974237762
